# The workbook was reopened from a new project folder (teampartition) and
# the lone data sheet was renamed from the default "Sheet1" to "2025" to
# reflect the season of player data it now holds.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "2025"
